# Auto-generated Excel COM-interop script to apply cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.527.53"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "2.484.15"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.86%  "
$ws.Range("E12").Value = "  +3.19%  "
$ws.Range("D13").Value = "2.866.05"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +10.41%  "
$ws.Range("D16").Value = "2.473.70"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").Value = "41.556.00"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("D20").Value = "0.0₃0935"
$ws.Range("E20").Value = "  +2.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("E24").Value = "  -0.45%  "
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.01%  "
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("E29").Value = "  +1.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.11%  "
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("E33").Value = "  +1.45%  "
$ws.Range("E34").Value = "  +3.08%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.07%  "
$ws.Range("B36").Value = "ApeXProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.104"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.70%  "
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "1.982.96"
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.11%  "
$ws.Range("D48").Value = "2.722.78"
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "98.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.09%  "
